$wb = $excel.ActiveWorkbook

# sheet -> list of (row, col, value_or_null, action)
# action: "set" -> set the value (numeric); "clear" -> clear the cell entirely (remove)

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 1600          # H32
$ws.Cells.Item(32, 10).Value = 0            # J32
$ws.Cells.Item(32, 12).Value = 0            # L32
$ws.Cells.Item(32, 14).ClearContents()      # N32 removed

$ws.Cells.Item(53, 8).Value = 729.7727      # H53
$ws.Cells.Item(53, 9).Value = 681.0769      # I53
$ws.Cells.Item(53, 10).Value = 800.1111     # J53
$ws.Cells.Item(53, 11).Value = 681.0769     # K53
$ws.Cells.Item(53, 12).Value = 800.1111     # L53
$ws.Cells.Item(53, 13).Value = -44.07690000000002  # M53
$ws.Cells.Item(53, 14).Value = -2074.1111   # N53

$ws.Cells.Item(106, 8).Value = 4966.3335    # H106
$ws.Cells.Item(106, 9).Value = 4966.3335    # I106
$ws.Cells.Item(106, 11).Value = 4966.3335   # K106
$ws.Cells.Item(106, 13).Value = -4335.3335  # M106

$ws.Cells.Item(107, 8).Value = 212.625      # H107
$ws.Cells.Item(107, 9).Value = 225.14285    # I107
$ws.Cells.Item(107, 11).Value = 225.14285   # K107
$ws.Cells.Item(107, 13).Value = 1694.85715  # M107

$ws.Cells.Item(132, 8).Value = 1929.5294    # H132
$ws.Cells.Item(132, 9).Value = 1414.4286    # I132
$ws.Cells.Item(132, 10).Value = 4333.3335   # J132
$ws.Cells.Item(132, 11).Value = 4243.2858   # K132
$ws.Cells.Item(132, 12).Value = 13000.0005  # L132
$ws.Cells.Item(132, 13).Value = -1713.2858  # M132
$ws.Cells.Item(132, 14).Value = -18060.0005 # N132

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 1051   # H5
$ws.Cells.Item(5, 9).Value = 1000   # I5
$ws.Cells.Item(5, 10).Value = 1102  # J5
$ws.Cells.Item(5, 11).Value = 1000  # K5
$ws.Cells.Item(5, 12).Value = 1102  # L5
$ws.Cells.Item(5, 13).Value = -888  # M5
$ws.Cells.Item(5, 14).Value = -1326 # N5

$ws.Cells.Item(32, 8).Value = 4114.353   # H32
$ws.Cells.Item(32, 9).Value = 4114.353   # I32
$ws.Cells.Item(32, 11).Value = 4114.353  # K32
$ws.Cells.Item(32, 13).Value = -3827.353 # M32

$ws.Cells.Item(107, 8).Value = 47497.5   # H107
$ws.Cells.Item(107, 10).Value = 47497.5  # J107
$ws.Cells.Item(107, 12).Value = 47497.5  # L107
$ws.Cells.Item(107, 14).Value = -55177.5 # N107

$ws.Cells.Item(110, 8).Value = 5661.3335   # H110
$ws.Cells.Item(110, 9).Value = 5661.3335   # I110
$ws.Cells.Item(110, 11).Value = 5661.3335  # K110
$ws.Cells.Item(110, 13).Value = -3616.3335 # M110

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 1051   # H4
$ws.Cells.Item(4, 9).Value = 1000   # I4
$ws.Cells.Item(4, 10).Value = 1102  # J4
$ws.Cells.Item(4, 11).Value = 1000  # K4
$ws.Cells.Item(4, 12).Value = 1102  # L4
$ws.Cells.Item(4, 13).Value = -885  # M4
$ws.Cells.Item(4, 14).Value = -1332 # N4

$ws.Cells.Item(32, 8).Value = 25000    # H32
$ws.Cells.Item(32, 10).Value = 25000   # J32
$ws.Cells.Item(32, 12).Value = 25000   # L32
$ws.Cells.Item(32, 14).Value = -25768  # N32

$ws.Cells.Item(134, 8).Value = 6275.5713    # H134
$ws.Cells.Item(134, 9).Value = 5976.3335    # I134
$ws.Cells.Item(134, 11).Value = 17929.0005  # K134
$ws.Cells.Item(134, 13).Value = -15394.0005 # M134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 4894.8335   # H58
$ws.Cells.Item(58, 10).Value = 10734.5    # J58
$ws.Cells.Item(58, 12).Value = 10734.5    # L58
$ws.Cells.Item(58, 14).Value = -11140.5   # N58

$ws.Cells.Item(86, 8).Value = 11620000  # H86
$ws.Cells.Item(86, 10).Value = 2880     # J86
$ws.Cells.Item(86, 12).Value = 2880     # L86
$ws.Cells.Item(86, 14).Value = -5126    # N86

$ws.Cells.Item(89, 8).Value = 11620000  # H89
$ws.Cells.Item(89, 10).Value = 2880     # J89
$ws.Cells.Item(89, 12).Value = 14400    # L89
$ws.Cells.Item(89, 14).Value = -25632   # N89

$ws.Cells.Item(94, 8).Value = 2578.3      # H94
$ws.Cells.Item(94, 9).Value = 2772.3333   # I94
$ws.Cells.Item(94, 10).Value = 2287.25    # J94
$ws.Cells.Item(94, 11).Value = 2772.3333  # K94
$ws.Cells.Item(94, 12).Value = 2287.25    # L94
$ws.Cells.Item(94, 13).Value = -2321.3333 # M94
$ws.Cells.Item(94, 14).Value = -3189.25   # N94

$ws.Cells.Item(99, 8).Value = 4473.533    # H99
$ws.Cells.Item(99, 9).Value = 3787.375    # I99
$ws.Cells.Item(99, 10).Value = 5257.7144  # J99
$ws.Cells.Item(99, 11).Value = 3787.375   # K99
$ws.Cells.Item(99, 12).Value = 5257.7144  # L99
$ws.Cells.Item(99, 13).Value = -2289.375  # M99
$ws.Cells.Item(99, 14).Value = -8253.714400000001 # N99

$ws.Cells.Item(105, 8).Value = 1848.875  # H105
$ws.Cells.Item(105, 9).Value = 897.3333  # I105
$ws.Cells.Item(105, 10).Value = 2419.8   # J105
$ws.Cells.Item(105, 11).Value = 897.3333 # K105
$ws.Cells.Item(105, 12).Value = 2419.8   # L105
$ws.Cells.Item(105, 13).Value = 849.6667 # M105
$ws.Cells.Item(105, 14).Value = -5913.8  # N105

$ws.Cells.Item(112, 8).Value = 55000   # H112
$ws.Cells.Item(112, 10).Value = 55000  # J112
$ws.Cells.Item(112, 12).Value = 55000  # L112
$ws.Cells.Item(112, 14).Value = -57954 # N112 (new)

$ws.Cells.Item(126, 8).Value = 4473.533     # H126
$ws.Cells.Item(126, 9).Value = 3787.375     # I126
$ws.Cells.Item(126, 10).Value = 5257.7144   # J126
$ws.Cells.Item(126, 11).Value = 11362.125   # K126
$ws.Cells.Item(126, 12).Value = 15773.1432  # L126
$ws.Cells.Item(126, 13).Value = -8892.125   # M126
$ws.Cells.Item(126, 14).Value = -20713.1432 # N126

$ws.Cells.Item(136, 8).Value = 4894.8335  # H136
$ws.Cells.Item(136, 10).Value = 10734.5   # J136
$ws.Cells.Item(136, 12).Value = 32203.5   # L136
$ws.Cells.Item(136, 14).Value = -37303.5  # N136

$ws.Cells.Item(141, 8).Value = 392310.38   # H141
$ws.Cells.Item(141, 10).Value = 392310.38  # J141
$ws.Cells.Item(141, 12).Value = 392310.38  # L141
$ws.Cells.Item(141, 14).Value = -402670.38 # N141

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 0     # H17
$ws.Cells.Item(17, 10).Value = 0    # J17
$ws.Cells.Item(17, 12).Value = 0    # L17
$ws.Cells.Item(17, 14).ClearContents()  # N17 removed

$ws.Cells.Item(23, 8).Value = 699.8889    # H23
$ws.Cells.Item(23, 10).Value = 779.2857   # J23
$ws.Cells.Item(23, 12).Value = 2337.8571  # L23
$ws.Cells.Item(23, 14).Value = -2807.8571 # N23

$ws.Cells.Item(61, 8).Value = 296.33334         # H61
$ws.Cells.Item(61, 9).Value = 95.666664         # I61
$ws.Cells.Item(61, 10).Value = 697.6667         # J61
$ws.Cells.Item(61, 11).Value = 286.999992       # K61
$ws.Cells.Item(61, 12).Value = 2093.0001        # L61
$ws.Cells.Item(61, 13).Value = -71.99999200000002 # M61
$ws.Cells.Item(61, 14).Value = -2523.0001       # N61

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 863.7059   # H102
$ws.Cells.Item(102, 9).Value = 667.6875   # I102
$ws.Cells.Item(102, 11).Value = 667.6875  # K102
$ws.Cells.Item(102, 13).Value = 954.3125  # M102

$ws.Cells.Item(106, 8).Value = 28000    # H106
$ws.Cells.Item(106, 10).Value = 28000   # J106
$ws.Cells.Item(106, 12).Value = 28000   # L106
$ws.Cells.Item(106, 14).Value = -30524  # N106 (new)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1866.8462   # H22
$ws.Cells.Item(22, 9).Value = 1085.4445   # I22
$ws.Cells.Item(22, 10).Value = 3625       # J22
$ws.Cells.Item(22, 11).Value = 1085.4445  # K22
$ws.Cells.Item(22, 12).Value = 3625       # L22
$ws.Cells.Item(22, 13).Value = -790.4445000000001 # M22
$ws.Cells.Item(22, 14).Value = -4215      # N22

$ws.Cells.Item(27, 8).Value = 1866.8462   # H27
$ws.Cells.Item(27, 9).Value = 1085.4445   # I27
$ws.Cells.Item(27, 10).Value = 3625       # J27
$ws.Cells.Item(27, 11).Value = 1085.4445  # K27
$ws.Cells.Item(27, 12).Value = 3625       # L27
$ws.Cells.Item(27, 13).Value = -978.4445000000001 # M27
$ws.Cells.Item(27, 14).Value = -3839      # N27

$ws.Cells.Item(29, 8).Value = 0     # H29
$ws.Cells.Item(29, 10).Value = 0    # J29
$ws.Cells.Item(29, 12).Value = 0    # L29
$ws.Cells.Item(29, 14).ClearContents()  # N29 removed

$ws.Cells.Item(46, 8).Value = 3055.2778   # H46
$ws.Cells.Item(46, 9).Value = 2375        # I46
$ws.Cells.Item(46, 10).Value = 4415.8335  # J46
$ws.Cells.Item(46, 11).Value = 2375       # K46
$ws.Cells.Item(46, 12).Value = 4415.8335  # L46
$ws.Cells.Item(46, 13).Value = -2187      # M46
$ws.Cells.Item(46, 14).Value = -4791.8335 # N46

$ws.Cells.Item(93, 8).Value = 1500   # H93
$ws.Cells.Item(93, 9).Value = 1500   # I93
$ws.Cells.Item(93, 11).Value = 1500  # K93
$ws.Cells.Item(93, 13).Value = -252  # M93 (new)

$ws.Cells.Item(100, 8).Value = 1999   # H100
$ws.Cells.Item(100, 9).Value = 1999   # I100
$ws.Cells.Item(100, 11).Value = 1999  # K100
$ws.Cells.Item(100, 13).Value = -1458 # M100

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1150.125    # H126
$ws.Cells.Item(126, 9).Value = 903.9286    # I126
$ws.Cells.Item(126, 11).Value = 2711.7858  # K126
$ws.Cells.Item(126, 13).Value = -241.7857999999997 # M126

$ws.Cells.Item(132, 8).Value = 2598.4285   # H132
$ws.Cells.Item(132, 9).Value = 2237.8      # I132
$ws.Cells.Item(132, 11).Value = 6713.400000000001  # K132
$ws.Cells.Item(132, 13).Value = -4183.400000000001 # M132

$ws.Cells.Item(136, 8).Value = 2254.1428   # H136
$ws.Cells.Item(136, 9).Value = 2043        # I136
$ws.Cells.Item(136, 11).Value = 6129       # K136
$ws.Cells.Item(136, 13).Value = -3579      # M136
